# 20170715 - Commit and Push
# Insert the new "12 July 2017" lotto draw as the most recent result,
# pushing all the existing draw rows down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row just above the current most-recent draw (row 6),
# shifting every existing data row down by one.
$ws.Rows("6:6").Insert()

# Populate the newly inserted row with the 12 July 2017 draw.
$ws.Cells.Item(6, 1).Value = "12 July 2017"
$ws.Cells.Item(6, 2).Value = 1
$ws.Cells.Item(6, 3).Value = 24
$ws.Cells.Item(6, 4).Value = 31
$ws.Cells.Item(6, 5).Value = 33
$ws.Cells.Item(6, 6).Value = 35
$ws.Cells.Item(6, 7).Value = 40
$ws.Cells.Item(6, 8).Value = 4
